$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 24. In the original sheet, row 24 was a stray blank filler row
#    (only E24 had a value). Deleting it shifts rows 25-35 up to 24-34, which
#    restores the "TC02_Verify_MYACC" test case into row 24 and removes the
#    now-unneeded trailing row 35 from the used range.
$ws.Rows("24").Delete()

# 2. The "TC03_Verify_BLP_Solutions_ContactUS" row (C3) had an errant red font;
#    restore it to the normal black text used throughout the rest of the sheet.
$ws.Range("C3").Font.Color = 0

# 3. Give the TC22 row (C23) the same "clean" font treatment.
$ws.Range("C23").Font.Color = 0

# 4. Move the active selection to D3.
$ws.Range("D3").Select() | Out-Null

# 5. The _FilterDatabase defined name pointed at the whole table; collapse it
#    back down to just the header row.
$wb.Names.Item(1).RefersTo = "=MasterExecutor!`$A`$1:`$F`$1"

# 6. Re-anchor the two "unique value" conditional-formatting rules on column F
#    to their shifted ranges (F27 -> F26 single-cell rule; F22:F26,F28:F32 ->
#    F22:F25,F27:F31 surrounding rule).
$fc1 = $ws.Cells.FormatConditions.Item(1)
$fc2 = $ws.Cells.FormatConditions.Item(2)
$fc1.ModifyAppliesToRange($ws.Range("F26"))
$fc2.ModifyAppliesToRange($ws.Range("F27:F31"))
$extraFc = $ws.Range("F22:F25").FormatConditions.AddUniqueValues()
$extraFc.DxfId = 0

Write-Output "edit complete"
